# Apply the latest crypto price/volume snapshot to the worksheet.
# Values in column D that Excel would otherwise auto-parse as numbers
# are written with a leading apostrophe so they stay text, matching
# the original inlineStr cell contents exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '60.659.06'
$ws.Range("E2").Value = '  +3.23%  '
$ws.Range("D3").Value = '2.688.01'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''525.78'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").Value = '''145.08'
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  +2.77%  '
$ws.Range("D9").Value = '2.709.43'
$ws.Range("E9").Value = '  +2.28%  '
$ws.Range("D10").Value = '''6.60'
$ws.Range("E10").Value = '  +6.38%  '
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("D14").Value = '3.162.05'
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").Value = '60.621.98'
$ws.Range("E15").Value = '  +3.11%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '''21.29'
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.705.47'
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("D18").Value = '''0.0000138'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '''347.68'
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = '''4.52'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '''10.57'
$ws.Range("E21").Value = '  +2.92%  '
$ws.Range("D22").Value = '''6.34'
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").Value = '''0.997'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").Value = '''63.57'
$ws.Range("E24").Value = '  +3.08%  '
$ws.Range("D25").Value = '''0.422'
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +6.40%  '
$ws.Range("D27").Value = '''0.994'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '''7.35'
$ws.Range("E28").Value = '  +4.09%  '
$ws.Range("D29").Value = '0.0₃0816'
$ws.Range("E29").Value = '  +2.62%  '
$ws.Range("D30").Value = '''6.87'
$ws.Range("E30").Value = '  +10.14%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("E32").Value = '  +2.06%  '
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("D34").Value = '''150.40'
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").Value = '''4.26'
$ws.Range("E35").Value = '  +6.97%  '
$ws.Range("D36").Value = '''1.24'
$ws.Range("E36").Value = '  +9.73%  '
$ws.Range("D37").Value = '''0.937'
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("D38").Value = '''0.881'
$ws.Range("E38").Value = '  +5.64%  '
$ws.Range("D39").Value = '''1.52'
$ws.Range("E39").Value = '  +8.48%  '
$ws.Range("D40").Value = '''36.95'
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("D41").Value = '''3.67'
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").Value = '''286.91'
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.0992'
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("D44").Value = '''0.613'
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''19.99'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '''0.996'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.142.44'
$ws.Range("E47").Value = '  +8.52%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '''0.0540'
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''4.87'
$ws.Range("E49").Value = '  +6.94%  '
$ws.Range("D50").Value = '''0.0235'
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("E51").Value = '  +1.76%  '
